$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header names (SamplePortion/Result and WaitingTime/Temperature swap) ---
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"
$ws.Range("J1").Value = "WaitingTime"
$ws.Range("K1").Value = "Temperature"

# --- Row 2: type row gains unit annotations ---
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"
$ws.Range("K2").Value = "#integer,  unit:celsius"

# --- Row 3: new description / enum row ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#TempsAttenteTempsAttente"
$ws.Range("K3").Value = "#Desc:Temperature"
